# Updated cryptos list on Sat Jan 27 23:32:36 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for the coinranking rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").Value = "42.123.84"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "2.266.07"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "150.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +14,940.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "94.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.04%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  +0.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.09%  "

$ws.Range("E12").Value = "  +0.41%  "

$ws.Range("E13").Value = "  -1.93%  "

$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").Value = "2.614.49"
$ws.Range("E15").Value = "  -0.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.59%  "

$ws.Range("D17").Value = "2.272.35"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("E18").Value = "  +3.62%  "

$ws.Range("D19").Value = "41.990.66"
$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.69%  "

$ws.Range("E21").Value = "  +1.13%  "

$ws.Range("E22").Value = "  +0.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.13"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "244.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.72%  "

$ws.Range("E26").Value = "  +2.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("E29").Value = "  +0.64%  "

$ws.Range("E30").Value = "  +0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("E33").Value = "  +3.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.97%  "

$ws.Range("E38").Value = "  -0.87%  "

$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("E40").Value = "  +0.67%  "

$ws.Range("E41").Value = "  -1.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.40%  "

$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").Value = "2.006.36"

$ws.Range("E45").Value = "  +11.09%  "

$ws.Range("E46").Value = "  +1.24%  "

$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
